$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price scrape: D-column prices and E-column 1h-volume deltas refresh
# each run. Price strings that look like plain numbers (e.g. "0.9997") must be
# written as literal text -- like the upstream openpyxl writer does -- instead of
# being auto-parsed into floating-point numbers by Excel, so force Text format
# before assignment and then clear the formatting override back off afterwards
# (ClearFormats keeps the cells stored text value/type, only the style resets).
function Set-TextCell($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

$ws.Range('D2').Value = '30.274.46'
$ws.Range('D3').Value = '1.930.60'
$ws.Range('E3').Value = '  -0.35%  '
Set-TextCell 'D4' '0.9997'
$ws.Range('E4').Value = '  -0.10%  '
Set-TextCell 'D5' '249.07'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -1.01%  '
Set-TextCell 'D7' '0.9995'
$ws.Range('E7').Value = '  -0.05%  '
Set-TextCell 'D8' '0.3201'
$ws.Range('E8').Value = '  -4.14%  '
$ws.Range('E9').Value = '  -2.68%  '
Set-TextCell 'D10' '0.07105'
$ws.Range('E10').Value = '  -4.08%  '
Set-TextCell 'D11' '0.7921'
$ws.Range('E11').Value = '  -2.82%  '
Set-TextCell 'D12' '0.07989'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '1.925.11'
$ws.Range('E13').Value = '  -0.62%  '
Set-TextCell 'D14' '5.393'
$ws.Range('E14').Value = '  -2.25%  '
Set-TextCell 'D15' '94.83'
$ws.Range('E15').Value = '  -0.38%  '
Set-TextCell 'D16' '14.68'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').Value = '30.254.52'
$ws.Range('E17').Value = '  -0.25%  '
Set-TextCell 'D18' '256.25'
$ws.Range('E18').Value = '  +0.60%  '
Set-TextCell 'D19' '0.000008038'
$ws.Range('E19').Value = '  -3.84%  '
Set-TextCell 'D20' '5.775'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('D21').Value = '2.176.72'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('E22').Value = '  -0.05%  '
Set-TextCell 'D23' '0.9994'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  -1.89%  '
Set-TextCell 'D25' '9.541'
$ws.Range('E25').Value = '  -2.72%  '
Set-TextCell 'D26' '165.00'
$ws.Range('E26').Value = '  +2.81%  '
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('E28').Value = '  -6.51%  '
$ws.Range('E29').Value = '  -4.59%  '
$ws.Range('E30').Value = '  +0.83%  '
Set-TextCell 'D31' '1.528'
$ws.Range('E31').Value = '  -2.10%  '
Set-TextCell 'D32' '4.394'
$ws.Range('E32').Value = '  -1.17%  '
Set-TextCell 'D33' '4.132'
$ws.Range('E33').Value = '  -2.72%  '
Set-TextCell 'D34' '0.05140'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('E35').Value = '  +1.26%  '
Set-TextCell 'D36' '0.7453'
$ws.Range('E36').Value = '  -0.80%  '
Set-TextCell 'D37' '2.772'
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range('E38').Value = '  -2.06%  '
Set-TextCell 'D39' '2.796'
$ws.Range('E39').Value = '  -1.62%  '
Set-TextCell 'D40' '78.71'
$ws.Range('E40').Value = '  -1.21%  '
Set-TextCell 'D41' '6.359'
$ws.Range('E41').Value = '  -4.72%  '
Set-TextCell 'D42' '0.4511'
$ws.Range('E42').Value = '  -0.80%  '
Set-TextCell 'D43' '1.989'
$ws.Range('E43').Value = '  -1.92%  '
Set-TextCell 'D44' '0.8478'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('E45').Value = '  -0.07%  '
Set-TextCell 'D46' '100.48'
$ws.Range('E46').Value = '  -2.15%  '
Set-TextCell 'D47' '9.795'
$ws.Range('E47').Value = '  -0.44%  '
Set-TextCell 'D48' '7.428'
$ws.Range('E48').Value = '  +0.44%  '
Set-TextCell 'D49' '36.66'
$ws.Range('E49').Value = '  -0.52%  '
Set-TextCell 'D50' '948.94'
$ws.Range('E50').Value = '  +10.07%  '
Set-TextCell 'D51' '0.06103'
$ws.Range('E51').Value = '  +1.82%  '
